$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Pruebas Realizadas" (column B) values for newly reported days
$ws.Range("B37").Value = 451
$ws.Range("B38").Value = 549
$ws.Range("B39").Value = 598
$ws.Range("B40").Value = 615
$ws.Range("B43").Value = 838

# Update the active selection to match the latest edited cell
$ws.Range("D38").Select()
